# Move the "For more info:" slide (position 2) to the end of the deck
# (after "Conclusions"), shifting every other slide up by one position.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$s.MoveTo($p.Slides.Count)
